$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (A1:J1)
$headers = @('province', ' Average income in the province', 'Population', ' Hazard (protection)', 'Exposure, poor people', 'Exposure, non-poor people', 'Asset vulnerability (poor people)', 'Asset vulnerability (non-poor people)', 'Socio-economic capacity', 'Risk to well-being')
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# I1:J1 are new header cells beyond the original A1:H1 range - copy the
# existing bold/border/centered header style from H1 onto them.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @(
    @('Abra', 133.688, 240135.244121424, 10.57, 22.2709275205961, 17.6335096604643, 15.2463285068896, 10, 75.1829916026878, 1.55646412528839),
    @('Agusan Del Norte', 179.014, 661728.454375173, 9.41, 37.1583320246628, 29.7393040328347, 36.3700034477779, 10, 79.7353503966536, 3.30685847178952),
    @('Agusan Del Sur', 126.492, 677779.68215436, 8.61, 13.9710440043128, 10.4613024935837, 29.0165593367945, 10, 56.1630386775849, 2.07793020290394),
    @('Apayao', 158.732, 116022.676390039, 6.51, 0.46, 0.304693688066755, 13.1939161763218, 10, 71.2496732637461, 0.0486244079857629),
    @('Batangas', 219.272, 2485002.05546474, 7.73, 0.2645, 0.221888963667484, 23.0356790840862, 10, 143.221300306026, 0.0133920570990175),
    @('Bukidnon', 145.749, 1353046.84430757, 9.82, 2.93531071612924, 2.1845799395281, 22.1694207791634, 10, 63.2524027580958, 0.291115179156782),
    @('Bulacan', 253.279, 3086241.84919508, 6.51, 1.24384889842907, 1.06881492470657, 41.8673134826198, 10.4395897258309, 191.226776325474, 0.0586168150176135),
    @('Cagayan', 179.424, 1153021.19591091, 9.34, 13.1684530678724, 11.029222045276, 29.4633115961937, 10, 101.874554000675, 0.819512727560968),
    @('Camarines Sur', 152.074, 1881962.89125995, 15.36, 12.2914883091891, 9.5666083772646, 33.8661767479629, 10, 67.1899542930178, 0.836887591655869),
    @('Capiz', 177.386, 733558.367597268, 15.94, 15.5045809179917, 12.7049561859472, 37.536961186832, 13.063852224195, 96.5409522283178, 0.801320978702004),
    @('Davao Del Norte', 156.283, 992308.737024673, 7.9, 2.40547700159559, 1.93413830291244, 30.6435266609834, 10, 72.0347924187506, 0.274576543988225),
    @('Ifugao', 149.062, 197584.190285786, 10.64, 8.21312668449679, 6.35173736056213, 14.6174762523903, 10, 81.4628496831763, 0.512846863630275),
    @('Ilocos Norte', 192.382, 579429.044110523, 48.14, 9.15433458917778, 7.8293585360055, 37.6494276892936, 11.069009836215, 138.007218501744, 0.0872163307808523),
    @('Ilocos Sur', 175.91, 672277.192943119, 27.75, 5.69738462250231, 4.79890201058604, 31.3748032009296, 10, 110.689598385586, 0.110033647286979),
    @('Isabela', 169.871, 1533715.59215201, 10.29, 11.3951572623664, 9.42989583353871, 25.7401509648132, 10, 93.2031154427546, 0.703163175058067),
    @('Kalinga', 170.385, 207635.165378523, 7.5, 7.93255398032548, 6.51967692084141, 17.3222304410523, 10, 106.521624349302, 0.550975509022513),
    @('La Union', 206.399, 759942.515079343, 23.79, 2.2715612175992, 1.90802973192671, 29.9801918728384, 10, 131.608113378751, 0.0422089228075164),
    @('Laguna', 249.092, 2838393.48641403, 8.78, 2.23513682475805, 1.92368476980292, 40.0677709602969, 10.5971668934738, 195.424997209117, 0.077007797950321),
    @('Lanao Del Norte', 164.672, 969716.979767905, 8.09, 1.26697409880717, 0.942551704195964, 22.1544883688664, 10, 70.2784935504147, 0.133309798085049),
    @('Lanao Del Sur', 131.981, 962426.556333156, 10.34, 5.819, 2.92559619824276, 24.4783044022428, 10, 39.6643490328256, 0.922094080716499),
    @('Maguindanao', 140.218, 1274192.94879532, 7.23, 9.26374475357312, 6.32247566275061, 28.127287626786, 10, 47.8597069688333, 1.87910204241408),
    @('Negros Occidental', 145.959, 2981577.20469059, 14.13, 4.66910471243158, 3.76909162168267, 35.7654244532121, 12.2579410501458, 75.1327372787753, 0.348581340707299),
    @('North Cotabato', 136.889, 1288466.35114615, 8.5, 11.8550672732701, 8.60950902519667, 28.8353730420083, 10, 55.5039011464904, 1.77736614293347),
    @('Nueva Ecija', 169.592, 2020505.5211558, 9.92, 18.7151631670846, 15.4518327134671, 20.8612015839051, 10, 101.04032379541, 1.06611049569636),
    @('Nueva Vizcaya', 214.169, 433165.204738251, 14.25, 23.630374211603, 19.7450352096992, 28.5573966207095, 10, 123.330849624189, 0.779238761029491),
    @('Pampanga', 235.621, 2444446.17134412, 8.49, 13.1917533823871, 11.3299454972311, 41.2615084591602, 10.3966925918505, 174.236608337867, 0.523405744510119),
    @('Pangasinan', 182.388, 2854701.4824351, 17.58, 18.3943328843056, 15.3812539629078, 28.1410282253697, 10, 111.91154184398, 0.549739219333778),
    @('Quezon', 179.092, 2055103.11595472, 17.61, 1.5079997366952, 1.23653705180482, 19.0092306563477, 10, 106.821592180453, 0.0449670311961865),
    @('Rizal', 294.404, 2678553.46348061, 16.23, 2.87917767015997, 2.47924740891148, 40.7993142799703, 10.6425721398091, 243.898686506192, 0.0428811436081932),
    @('Samar', 134.968, 753363.084895201, 11.87, 9.26762422927215, 6.84867302200992, 29.0019172186821, 10, 59.1009483492337, 0.935019465712449),
    @('Sarangani', 107.685, 518718.742502788, 23.94, 0.0575, 0.0414913493520872, 28.5538860754928, 10, 46.925854195331, 0.00396810518502108),
    @('South Cotabato', 191.360000000001, 1424914.17073343, 6.51, 1.15588206417351, 0.934225826639033, 34.9683979732136, 12.7583396548576, 93.4470460179474, 0.146011279970418),
    @('Sultan Kudarat', 113.736, 784135.419521753, 7.14, 8.9449715715031, 6.68059193122331, 30.343632332348, 10, 49.7503231676702, 1.93058097801387),
    @('Tarlac', 190.392, 1318604.09835682, 10.72, 14.6354343059504, 12.3476788908042, 26.52806518946, 10, 122.727110766257, 0.63729268341886),
    @('Zamboanga Del Sur', 168.131, 1841096.51798489, 15.31, 5.5860028179362, 4.51406553503674, 40.3643010697264, 12.9130326981283, 78.3911350579756, 0.384079334147517),
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
